$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "21/03/2023"

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 352
$ws.Range("D2").Value = 353
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 365.3
$ws.Range("J2").Value = -3.367095537914033

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 173
$ws.Range("D3").Value = 173
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 176
$ws.Range("J3").Value = -1.704545454545459

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 133.3333333333333

$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 131
$ws.Range("D5").Value = 154
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 67
$ws.Range("J5").Value = 129.8507462686567

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 28
$ws.Range("D6").Value = 28
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 32
$ws.Range("J6").Value = -12.5

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 31
$ws.Range("D7").Value = 34
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 49
$ws.Range("J7").Value = -30.61224489795918

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 143
$ws.Range("D8").Value = 154
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 38
$ws.Range("J8").Value = 305.2631578947368

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 18
$ws.Range("D9").Value = 18
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 82
$ws.Range("J9").Value = -78.04878048780488

$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 146
$ws.Range("D10").Value = 179
$ws.Range("E10").Value = 28
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 297
$ws.Range("J10").Value = -39.73063973063973

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 144
$ws.Range("D11").Value = 145
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 94
$ws.Range("J11").Value = 54.25531914893617

$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 231
$ws.Range("D12").Value = 257
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 17
$ws.Range("I12").Value = 341.9
$ws.Range("J12").Value = -24.83182217022521

$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 265
$ws.Range("J13").Value = -95.47169811320755

$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 246
$ws.Range("D14").Value = 259
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 10
$ws.Range("I14").Value = 410
$ws.Range("J14").Value = -36.82926829268293

$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 93
$ws.Range("D15").Value = 94
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 107
$ws.Range("J15").Value = -12.14953271028038

$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 89
$ws.Range("D16").Value = 113
$ws.Range("E16").Value = 21
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 119
$ws.Range("J16").Value = -5.042016806722693

$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 38
$ws.Range("D17").Value = 38
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = -24

$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 33.33333333333333

$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 17
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 38
$ws.Range("J20").Value = -55.26315789473684
